$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)

# Replace the paragraph text (runs 1+2: "**ID__AFFARS_5304_topic_2__ID**" + " ")
# with the single new id, dropping the trailing space run entirely.
$r = $p.Range
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.Text = "**ID__AFFARS_SUBPART_5304_1__ID**"

# Update the paragraph indentation (w:ind w:left) from 120 -> 225 twips.
$p.Format.LeftIndent = 11.25

# Add a paragraph border (pBdr) on all four sides with 5-twip spacing.
$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
